$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 136, shifting existing rows 136-141 down to 137-142
$ws.Rows.Item(136).Insert()

# Populate the new row 136 with data
$ws.Range("A136").Value = 3
$ws.Range("B136").Value = "Femacal de La Calera"
$ws.Range("C136").Value = "Coquimbo"
$ws.Range("D136").Value = 44509
$ws.Range("E136").Value = 5
$ws.Range("F136").Value = 100112010
$ws.Range("G136").Value = "Achicoria"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 125
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 6200
$ws.Range("M136").Value = 6104
$ws.Range("N136").Value = "$/caja 16 unidades"
$ws.Range("O136").Value = "Provincia de Quillota"
$ws.Range("P136").Value = 382
$ws.Range("Q136").Value = 16
$ws.Range("R136").Value = "Hortaliza"
